$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5:D5").NumberFormat = "@"

$ws.Range("A5").Value = "118450"
$ws.Range("B5").Value = "1008617899"
$ws.Range("C5").Value = "17706586"
$ws.Range("D5").Value = "6020"
